# "index is ready" — 02-02-2024
#
# The IsReleased / is_released column headers are being renamed:
#   - English header  (row 1, col AN): "IsReleased"   -> "is_released"
#   - Arabic label     (row 2, col AN): "مفرج عنه؟"     -> "مفرج عنه"
# (the arrest_end_date / تاريخ الإفراج pair in column N keeps its text,
#  it just shifts shared-string slots as a side effect of the sst edit)
#
# Also nudge the saved view state to match the author's last on-screen
# position: scrolled a little further left and one row up from the
# previous selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("AN1").Value = "is_released"
$ws.Range("AN2").Value = "مفرج عنه"

# Move the selection the way the author's session ended up (N5 -> N4) and
# scroll the viewport back towards column B (topLeftCell E1 -> B1).
[void]$ws.Range("N4").Select()
try {
    $excel.ActiveWindow.ScrollColumn = 2
    $excel.ActiveWindow.ScrollRow = 1
} catch {
    # view-scroll position isn't always settable in every host; ignore if so
}
